$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.061582803726196
$ws.Range("B1").Value = 2.259390830993652
$ws.Range("C1").Value = 2.570086717605591
$ws.Range("D1").Value = 3.895766019821167
$ws.Range("E1").Value = 1.034889221191406
